$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 5.3
$ws.Range("Q2").Value = 1.81
$ws.Range("S2").Value = 3.15
$ws.Range("Z2").Value = 11.5
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AG2").Value = 22
$ws.Range("AI2").Value = 1000
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 3.7
$ws.Range("U4").Value = 2.26
$ws.Range("W4").Value = 1.91
$ws.Range("F5").Value = 1.91
$ws.Range("G5").Value = 1.99
$ws.Range("H5").Value = 3.95
$ws.Range("I5").Value = 4.5
$ws.Range("J5").Value = 3.65
$ws.Range("P5").Value = 1.74
$ws.Range("T5").Value = 1.88
$ws.Range("U5").Value = 1.91
$ws.Range("V5").Value = 1.28
$ws.Range("W5").Value = 2
$ws.Range("X5").Value = 15
$ws.Range("AB5").Value = 9.6
$ws.Range("AC5").Value = 9.6
$ws.Range("AG5").Value = 12.5
$ws.Range("I7").Value = 2.48
$ws.Range("N7").Value = 2.1
$ws.Range("P7").Value = 1.47
$ws.Range("Q7").Value = 2.2
$ws.Range("V7").Value = 1.67
$ws.Range("N8").Value = 1.84
$ws.Range("L11").Value = 1.32
$ws.Range("V11").Value = 1.95
$ws.Range("F12").Value = 3.9
$ws.Range("H12").Value = 1.92
$ws.Range("I12").Value = 2.04
$ws.Range("R12").Value = 1.59
$ws.Range("T12").Value = 1.48
$ws.Range("V12").Value = 1.97
$ws.Range("F13").Value = 2.12
$ws.Range("G13").Value = 2.18
$ws.Range("H13").Value = 3.75
$ws.Range("J13").Value = 3.6
$ws.Range("K13").Value = 3.65
$ws.Range("W13").Value = 1.84
$ws.Range("X13").Value = 16.5
$ws.Range("AF13").Value = 14
$ws.Range("AL13").Value = 32
$ws.Range("H14").Value = 1.93
$ws.Range("I14").Value = 2.44
$ws.Range("J14").Value = 2.68
$ws.Range("K14").Value = 6
$ws.Range("G15").Value = 2.54
$ws.Range("N15").Value = 1.31
$ws.Range("P15").Value = 1.3
$ws.Range("Q15").Value = 1.91
$ws.Range("T15").Value = 1.79
$ws.Range("U15").Value = 1.01
$ws.Range("W15").Value = 1.65
$ws.Range("T16").Value = 1.85
$ws.Range("U16").Value = 1.83
$ws.Range("N18").Value = 2.22
$ws.Range("Q19").Value = 1.5
$ws.Range("F20").Value = 2.8
$ws.Range("G20").Value = 3.25
$ws.Range("H20").Value = 2.54
$ws.Range("I20").Value = 3.2
$ws.Range("K20").Value = 4.1
$ws.Range("P20").Value = 1.74
$ws.Range("Q20").Value = 2.1
$ws.Range("V20").Value = 1.46
$ws.Range("W20").Value = 1.44
$ws.Range("K24").Value = 4.1
$ws.Range("M24").Value = 1.09
$ws.Range("N24").Value = 2.98
$ws.Range("Q24").Value = 2.16
$ws.Range("T24").Value = 1.9
$ws.Range("U24").Value = 1.9
$ws.Range("G25").Value = 3.15
$ws.Range("H25").Value = 2.8
$ws.Range("I25").Value = 3.65
$ws.Range("S25").Value = 4.3
$ws.Range("U25").Value = 1.78
$ws.Range("V25").Value = 1.43
$ws.Range("W25").Value = 1.47
$ws.Range("G26").Value = 3.65
$ws.Range("T26").Value = 1.87
$ws.Range("AA26").Value = 46
$ws.Range("AE26").Value = 34
$ws.Range("AJ26").Value = 75
$ws.Range("AK26").Value = 55
$ws.Range("AN26").Value = 60
$ws.Range("H27").Value = 6.4
$ws.Range("I27").Value = 8
$ws.Range("L27").Value = 1.34
$ws.Range("V27").Value = 1.12
$ws.Range("F28").Value = 2.66
$ws.Range("I28").Value = 3.25
$ws.Range("G29").Value = 2.88
$ws.Range("J29").Value = 3.4
$ws.Range("Q29").Value = 1.75
$ws.Range("T29").Value = 1.6
$ws.Range("U29").Value = 2
$ws.Range("V29").Value = 1.49
$ws.Range("W29").Value = 1.54
$ws.Range("G30").Value = 2.1
$ws.Range("K30").Value = 3.95
$ws.Range("M30").Value = 1.07
$ws.Range("U30").Value = 1.79
$ws.Range("V30").Value = 1.28
$ws.Range("W30").Value = 1.91
$ws.Range("AI30").Value = 75
$ws.Range("I32").Value = 5.7
$ws.Range("J32").Value = 2.78
$ws.Range("N32").Value = 2.06
$ws.Range("P32").Value = 2.06
$ws.Range("Q32").Value = 1.58
$ws.Range("F33").Value = 2.06
$ws.Range("G33").Value = 2.3
$ws.Range("J33").Value = 3.4
$ws.Range("K33").Value = 4.1
$ws.Range("O33").Value = 1.19
$ws.Range("P33").Value = 2.38
$ws.Range("T33").Value = 1.46
$ws.Range("U33").Value = 2.54
$ws.Range("W33").Value = 1.77
$ws.Range("AC33").Value = 9.4
$ws.Range("F35").Value = 2.16
$ws.Range("G35").Value = 2.3
$ws.Range("I35").Value = 3.35
$ws.Range("L35").Value = 1.23
$ws.Range("Q35").Value = 1.5
$ws.Range("S35").Value = 2.22
$ws.Range("V35").Value = 1.43
$ws.Range("W35").Value = 1.76
$ws.Range("N36").Value = 1.01
$ws.Range("P36").Value = 2.62
$ws.Range("S36").Value = 1.91
$ws.Range("T36").Value = 1.55
$ws.Range("U36").Value = 2.1
$ws.Range("P37").Value = 1.25
$ws.Range("S37").Value = 1.01
$ws.Range("G38").Value = 2.28
$ws.Range("H38").Value = 3
$ws.Range("I38").Value = 4.4
$ws.Range("Q38").Value = 1.12
$ws.Range("V38").Value = 1.29
$ws.Range("W38").Value = 1.78
$ws.Range("J39").Value = 4
$ws.Range("S39").Value = 2.22
$ws.Range("V39").Value = 1.41
$ws.Range("Q40").Value = 1.28
$ws.Range("T41").Value = 1.77
$ws.Range("U41").Value = 2.06
$ws.Range("G42").Value = 2.64
$ws.Range("K42").Value = 3.2
$ws.Range("S42").Value = 5.5
$ws.Range("U42").Value = 1.72
$ws.Range("W42").Value = 1.61
$ws.Range("F43").Value = 2
$ws.Range("G43").Value = 2.7
$ws.Range("J43").Value = 2.8
$ws.Range("K43").Value = 5.3
$ws.Range("L43").Value = 1.4
$ws.Range("V43").Value = 1.24
$ws.Range("W43").Value = 1.59
$ws.Range("H44").Value = 2.16
$ws.Range("J44").Value = 3.2
$ws.Range("L44").Value = 1.32
$ws.Range("N44").Value = 1.84
$ws.Range("P44").Value = 1.84
$ws.Range("N45").Value = 4.3
$ws.Range("P45").Value = 2.12
$ws.Range("T45").Value = 1.73
$ws.Range("U45").Value = 2.26
$ws.Range("AM45").Value = 85
$ws.Range("AO45").Value = 13.5
$ws.Range("F46").Value = 1.76
$ws.Range("G46").Value = 1.81
$ws.Range("V46").Value = 1.19
$ws.Range("W46").Value = 2.22
$ws.Range("P47").Value = 1.72
$ws.Range("T47").Value = 1.96
$ws.Range("X47").Value = 11
$ws.Range("I48").Value = 3.7
$ws.Range("L48").Value = 1.6
$ws.Range("N48").Value = 2.62
$ws.Range("O48").Value = 1.6
$ws.Range("Y48").Value = 9.4
$ws.Range("F49").Value = 2.86
$ws.Range("H49").Value = 2.22
$ws.Range("L49").Value = 1.33
$ws.Range("N49").Value = 1.8
$ws.Range("P49").Value = 1.79
$ws.Range("Q49").Value = 2
$ws.Range("F50").Value = 1.79
$ws.Range("G50").Value = 1.81
$ws.Range("H50").Value = 5.1
$ws.Range("AB50").Value = 8.800000000000001
$ws.Range("F51").Value = 1.9
$ws.Range("G51").Value = 2.36
$ws.Range("I51").Value = 5.2
$ws.Range("K51").Value = 5
$ws.Range("L51").Value = 1.37
$ws.Range("Q51").Value = 2.08
$ws.Range("R51").Value = 1.22
$ws.Range("S51").Value = 3.55
$ws.Range("V51").Value = 1.23
$ws.Range("W51").Value = 1.74
